$d = $word.ActiveDocument

# Locate the run containing "AS PS VS VV" and replace its text with "178047".
$range = $d.Content
$found = $range.Find.Execute("AS PS VS VV", $true, $false, $false, $false, $false, $true, 1, $false, "178047", 1)

if ($found) {
    # $range now covers the replaced text ("178047"). Remove bold formatting
    # that the original "AS PS VS VV" run carried, matching the target XML
    # which drops <w:b/> and <w:bCs/>.
    $range.Font.Bold = 0
    $range.Font.BoldBi = 0
}
